$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Delete the old row 2 ("合肥·2024运动新春动漫庆典（全ip）") -- rows below shift up,
    # and the previously-last row (old row 8) disappears, shrinking the used range to A1:I7.
    $ws.Rows.Item(2).Delete()

    # A handful of "想去人数" (F column) counts were bumped for the surviving rows.
    $ws.Range("F2").Value = 365
    $ws.Range("F3").Value = 1220
    $ws.Range("F4").Value = 1394
    $ws.Range("F6").Value = 6089

    # Column A is a plain sequence number (not a formula) so Delete() carried the
    # old numbering up with each row; restore it to the original 1..6 sequence.
    $ws.Range("A2").Value = 1
    $ws.Range("A3").Value = 2
    $ws.Range("A4").Value = 3
    $ws.Range("A5").Value = 4
    $ws.Range("A6").Value = 5
    $ws.Range("A7").Value = 6
}

# Leave the originally-active sheet selected.
$wb.Worksheets.Item("展览").Activate()
